$d = $word.ActiveDocument

# Locate the body-level paragraph that contains the "{newlinetester}" placeholder
# (the same literal text also appears inside a table cell further down in the
# document, so stop at the first match, which is the body paragraph we want).
$targetPara = $null
$targetIdx = 0
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "{newlinetester}") {
        $targetPara = $para
        $targetIdx = $i
        break
    }
}

# Insert two new paragraphs right after it. InsertParagraphAfter clones the
# paragraph formatting (style + run properties) of the paragraph it is called on,
# so the new paragraphs automatically come out bold with the "para4" style,
# matching the surrounding content. Re-fetch each paragraph from the Paragraphs
# collection by index after every mutation instead of reusing a stale Range.
$rng = $targetPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newPara1 = $d.Paragraphs($targetIdx + 1)
$rng2 = $newPara1.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

# The second inserted paragraph (now at $targetIdx + 2) gets the new placeholder
# text, "{ampersand}".
$newPara2 = $d.Paragraphs($targetIdx + 2)
$newPara2.Range.InsertBefore("{ampersand}")
